$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 20.655921
$ws.Range("N2").Value = 61.967763
$ws.Range("O2").Value = 0.3846075519180752
$ws.Range("P2").Value = 0.3846075519180752
$ws.Range("Q2").Value = 0.3669317806439999
$ws.Range("R2").Value = 3.302386025796
$ws.Range("S2").Value = 0.3846075519180752
$ws.Range("T2").Value = 0.3846075519180752

# Row 3 updates
$ws.Range("O3").Value = 0.3973836246290912
$ws.Range("P3").Value = 0.3973836246290912
$ws.Range("S3").Value = 0.3973836246290912
$ws.Range("T3").Value = 0.3973836246290912

# Row 4 updates
$ws.Range("M4").Value = 11.70848833333333
$ws.Range("N4").Value = 35.125465
$ws.Range("O4").Value = 0.2180088234528336
$ws.Range("P4").Value = 0.2180088234528336
$ws.Range("S4").Value = 0.2180088234528336
$ws.Range("T4").Value = 0.2180088234528336
